$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E28").Value = 181

$ws.Range("E44").Value = 279

$ws.Range("E45").Value = 125

$ws.Range("E46").Value = 281
$ws.Range("F46").Value = 151
$ws.Range("H46").Value = 151
